$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (Client Id, Candidate ID, User Name, Exam Password, First Name, Last Name)
$ws.Range("A2").Value = "OaJZz461"
$ws.Range("B2").Value = 231102307
$ws.Range("C2").Value = "bsemrgy77"
$ws.Range("D2").Value = 'MJc$q&75'
$ws.Range("F2").Value = "fsFlDTuI"
$ws.Range("G2").Value = "HASF"
